# Add the "step one" body paragraph to the sdd document and clean up
# the two heading runs' direct formatting, matching the authored edit.
#
# We build the full replacement XML for the document body (the two
# existing numbered headings plus the new descriptive paragraph in
# between) and push it in one shot via Range.InsertXML so the new
# paragraph is created with exactly the formatting it needs instead of
# inheriting the bold/blue heading look from the paragraph it would
# otherwise split off from.

$d = $word.ActiveDocument

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:eastAsia="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="091E42"/>
        </w:rPr>
        <w:t>Gather requirements and initial information</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" w:firstLine="720"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">This project </w:t>
      </w:r>
      <w:r>
        <w:t>requires</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> that the lead and the </w:t>
      </w:r>
      <w:r>
        <w:t>designer to make a driver for the rest of the team and work on the first two steps of the software design document.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>Additionally</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, the </w:t>
      </w:r>
      <w:r>
        <w:t>designer must</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> create the class diagrams for the project. The SWE will code the even number problems while the </w:t>
      </w:r>
      <w:r>
        <w:t>tester will</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>code</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> the </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">odd numbered problems. The SWE should write AAA test in </w:t>
      </w:r>
      <w:r>
        <w:t>English</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> for the problems the Testers codes. While the Tester writes AAA tests in English for the </w:t>
      </w:r>
      <w:r>
        <w:t>problems</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> the </w:t>
      </w:r>
      <w:r>
        <w:t>SWE coded.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Since only Java’s collections classes can be used </w:t>
      </w:r>
      <w:r>
        <w:t>when coding</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> these </w:t>
      </w:r>
      <w:r>
        <w:t>problems</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, the SWE and the Tester must code in </w:t>
      </w:r>
      <w:r>
        <w:t>Java</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> All team members should </w:t>
      </w:r>
      <w:r>
        <w:t>commit</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> and push their work to the repo, through issuing pull </w:t>
      </w:r>
      <w:r>
        <w:t>requests</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> that the lead is responsible for </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">verifying. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">The lead is Kylan, the Designer is Taylor, the SWE is DeAnna, and the Tester is Brandon. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI" w:eastAsia="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="091E42"/>
        </w:rPr>
        <w:t>Outline the system architecture</w:t>
      </w:r>
    </w:p>
    <w:sectPr>
      <w:pgSz w:w="12240" w:h="15840"/>
      <w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/>
      <w:cols w:space="720"/>
      <w:docGrid w:linePitch="360"/>
    </w:sectPr>
  </w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d.Content.InsertXML($xml)
